# 300748-金力永磁.xlsx — add 2022-Q4 data
#
# 1) "总计" (summary) sheet: insert a new row for "2022-Q4" right after the
#    header row, shifting all existing quarter rows down by one, and fix up
#    the sequential index in column A for the row that used to be last.
# 2) Insert a brand-new "2022-Q4" worksheet (with the same layout/styling as
#    the existing per-quarter fund sheets) positioned right before the
#    "2022-Q3" tab, and populate it with the Q4 fund holdings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet — insert the new 2022-Q4 row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push rows 2..8 down to 3..9, leaving row 2 blank.
$summary.Rows("2:2").Insert()

# Row 2 inherited the header's formatting via the default insert; copy the
# correct data-row formatting (from the row immediately below, which holds
# what used to be row 2 / "2022-Q3") onto the new blank row instead.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

# Fill in the new 2022-Q4 summary row.
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 8
$summary.Cells.Item(2,4).Value = 0.7

# Rows("2:2").Insert() shifted the old rows (and their column-A running
# index) down by one physical row along with the quarter data. In the
# target layout, column A is pinned to the row position (0,1,2,...) rather
# than travelling with the data, so restore rows 3..9's index accordingly
# (row 9 is the former last row, "2021-Q1", whose index continues to 7).
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(8,1).Value = 6
$summary.Cells.Item(9,1).Value = 7

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q4" fund-holdings sheet
# ---------------------------------------------------------------------
$sheets = $wb.Worksheets
$q3 = $sheets.Item("2022-Q3")

# Duplicate the "2022-Q3" sheet (same headers/column styling) placed right
# before it, then rename the copy and trim it down to the 8 data rows we
# actually need (the source sheet has many more fund rows).
$q3.Copy($q3)
$q4 = $sheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"
$q4.Rows("10:29").Delete()

# Force text storage for the numeric-looking columns that must stay text
# (fund codes, AUM, position %, weight %, and most "market value" cells),
# matching how the other quarter sheets store this data.
$q4.Range("B2:F9").NumberFormat = "@"
$q4.Range("G2:G6").NumberFormat = "@"

$data = @(
    @(0, "470009", "汇添富民营活力混合A",             "24.48", "93.48", "2.69", "0.6585", 8),
    @(1, "163110", "申万菱信量化小盘股票（LOF）A",     "5.10",  "92.35", "0.55", "0.0280", 9),
    @(2, "009128", "明亚价值长青混合A",               "0.39",  "52.21", "3.27", "0.0128", 2),
    @(3, "004918", "中银证券祥瑞混合C",               "0.05",  "76.82", "2.77", "0.0014", 4),
    @(4, "004917", "中银证券祥瑞混合A",               "0.03",  "76.82", "2.77", "0.0008", 4),
    @(5, "009129", "明亚价值长青混合C",               "0.00",  "52.21", "3.27", 0,        2),
    @(6, "960014", "汇添富民营活力混合 O",            "0.00",  "93.48", "2.69", 0,        8),
    @(7, "013918", "申万菱信量化小盘股票（LOF）C",     "0.00",  "92.35", "0.55", 0,        9)
)

$r = 2
foreach ($row in $data) {
    $q4.Cells.Item($r,1).Value = $row[0]
    $q4.Cells.Item($r,2).Value = $row[1]
    $q4.Cells.Item($r,3).Value = $row[2]
    $q4.Cells.Item($r,4).Value = $row[3]
    $q4.Cells.Item($r,5).Value = $row[4]
    $q4.Cells.Item($r,6).Value = $row[5]
    $q4.Cells.Item($r,7).Value = $row[6]
    $q4.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Restore the original active sheet ("总计")
# ---------------------------------------------------------------------
$summary.Activate()

Write-Output "2022-Q4 sheet + summary row added"
